$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.813.19"
$ws.Range("E2").Value = "  +0.18%  "
$ws.Range("D3").Value = "3.148.84"
$ws.Range("E3").Value = "  +1.49%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'531.72"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.69%  "
$ws.Range("D6").Value = "'140.69"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.16%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "3.144.98"
$ws.Range("E8").Value = "  +1.37%  "
$ws.Range("E9").Value = "  +2.83%  "
$ws.Range("D10").Value = "'7.22"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.72%  "
$ws.Range("E11").Value = "  +0.30%  "
$ws.Range("D12").Value = "'0.401"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +4.92%  "
$ws.Range("D13").Value = "3.688.49"
$ws.Range("E13").Value = "  +1.44%  "
$ws.Range("E14").Value = "  +2.97%  "
$ws.Range("D15").Value = "'25.70"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.73%  "
$ws.Range("E16").Value = "  +0.99%  "
$ws.Range("D17").Value = "57.916.68"
$ws.Range("E17").Value = "  +0.21%  "
$ws.Range("D18").Value = "3.159.90"
$ws.Range("E18").Value = "  +1.63%  "
$ws.Range("D19").Value = "'6.12"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.48%  "
$ws.Range("D20").Value = "'12.83"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.34%  "
$ws.Range("E21").Value = "  -0.22%  "
$ws.Range("D22").Value = "'355.11"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +5.37%  "
$ws.Range("D23").Value = "'0.998"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.06%  "
$ws.Range("D24").Value = "'68.53"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.15%  "
$ws.Range("E25").Value = "  +0.73%  "
$ws.Range("E26").Value = "  +1.45%  "
$ws.Range("E27").Value = "  +0.21%  "
$ws.Range("D28").Value = "0.0₃0939"
$ws.Range("E28").Value = "  +1.41%  "
$ws.Range("D29").Value = "'7.46"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.51%  "
$ws.Range("D30").Value = "'0.999"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.06%  "
$ws.Range("D31").Value = "'6.38"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.33%  "
$ws.Range("E32").Value = "  +2.28%  "
$ws.Range("E33").Value = "  +1.98%  "
$ws.Range("E34").Value = "  +1.31%  "
$ws.Range("D35").Value = "'4.90"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +6.21%  "
$ws.Range("D36").Value = "'157.51"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.15%  "
$ws.Range("D37").Value = "'6.21"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.39%  "
$ws.Range("D38").Value = "'26.28"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.38%  "
$ws.Range("E39").Value = "  -0.85%  "
$ws.Range("D40").Value = "'0.0673"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.04%  "
$ws.Range("E41").Value = "  +11.46%  "
$ws.Range("D42").Value = "'4.08"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +5.03%  "
$ws.Range("D43").Value = "'0.704"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.15%  "
$ws.Range("D44").Value = "3.187.74"
$ws.Range("E44").Value = "  +1.30%  "
$ws.Range("E45").Value = "  +6.17%  "
$ws.Range("E46").Value = "  -0.41%  "
$ws.Range("D47").Value = "'0.999"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.10%  "
$ws.Range("D48").Value = "2.333.01"
$ws.Range("D49").Value = "'1.00"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.33%  "
$ws.Range("E50").Value = "  +1.29%  "
$ws.Range("D51").Value = "'20.36"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.19%  "
